# Swap the order of "dnasr281@gmail.com" and "System" in the Recorded By
# column (G) wherever both appear together, e.g.
#   "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

# Determine the last used row in column G (Recorded By).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

$changed = 0
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
        $changed = $changed + 1
    }
}

Write-Host "Updated cells:" $changed
